$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2394516.2
$ws.Range("J17").Value = 2430794.2
$ws.Range("L17").Value = 7292382.600000001
$ws.Range("N17").Value = -7292718.600000001
$ws.Range("H31").Value = 99
$ws.Range("I31").Value = 99
$ws.Range("K31").Value = 297
$ws.Range("M31").Value = -67
$ws.Range("H34").Value = 4238.5
$ws.Range("I34").Value = 2266.2
$ws.Range("J34").Value = 14100
$ws.Range("K34").Value = 2266.2
$ws.Range("L34").Value = 14100
$ws.Range("M34").Value = -2063.2
$ws.Range("N34").Value = -14506
$ws.Range("H36").Value = 4238.5
$ws.Range("I36").Value = 2266.2
$ws.Range("J36").Value = 14100
$ws.Range("K36").Value = 2266.2
$ws.Range("L36").Value = 14100
$ws.Range("M36").Value = -1551.2
$ws.Range("N36").Value = -15530
$ws.Range("H76").Value = 3774.3333
$ws.Range("I76").Value = 3263.9092
$ws.Range("J76").Value = 4125.25
$ws.Range("K76").Value = 3263.9092
$ws.Range("L76").Value = 4125.25
$ws.Range("M76").Value = -2948.9092
$ws.Range("N76").Value = -4755.25
$ws.Range("H79").Value = 3774.3333
$ws.Range("I79").Value = 3263.9092
$ws.Range("J79").Value = 4125.25
$ws.Range("K79").Value = 3263.9092
$ws.Range("L79").Value = 4125.25
$ws.Range("M79").Value = -2171.9092
$ws.Range("N79").Value = -6309.25
$ws.Range("H80").Value = 689.6
$ws.Range("I80").Value = 587.8461
$ws.Range("J80").Value = 878.5714
$ws.Range("K80").Value = 1763.5383
$ws.Range("L80").Value = 2635.7142
$ws.Range("M80").Value = -765.5382999999999
$ws.Range("N80").Value = -4631.7142
$ws.Range("H83").Value = 689.6
$ws.Range("I83").Value = 587.8461
$ws.Range("J83").Value = 878.5714
$ws.Range("K83").Value = 5290.6149
$ws.Range("L83").Value = 7907.1426
$ws.Range("M83").Value = -298.6148999999996
$ws.Range("N83").Value = -17891.1426
$ws.Range("H109").Value = 53500
$ws.Range("J109").Value = 53500
$ws.Range("L109").Value = 53500
$ws.Range("N109").Value = -56274
$ws.Range("H132").Value = 2364.2307
$ws.Range("I132").Value = 2364.2307
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7092.6921
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4562.6921
$ws.Range("N132").Value = $null
$ws.Range("H137").Value = 1873.2632
$ws.Range("I137").Value = 1852.8
$ws.Range("J137").Value = 1950
$ws.Range("K137").Value = 5558.4
$ws.Range("L137").Value = 5850
$ws.Range("M137").Value = -3008.4
$ws.Range("N137").Value = -10950

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21930
$ws.Range("J44").Value = 21930
$ws.Range("L44").Value = 21930
$ws.Range("N44").Value = -22906

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("H80").Value = 287.0476
$ws.Range("I80").Value = 48.333332
$ws.Range("J80").Value = 382.53333
$ws.Range("K80").Value = 48.333332
$ws.Range("L80").Value = 382.53333
$ws.Range("M80").Value = 949.666668
$ws.Range("N80").Value = -2378.53333
$ws.Range("H83").Value = 287.0476
$ws.Range("I83").Value = 48.333332
$ws.Range("J83").Value = 382.53333
$ws.Range("K83").Value = 241.66666
$ws.Range("L83").Value = 1912.66665
$ws.Range("M83").Value = 4750.33334
$ws.Range("N83").Value = -11896.66665
$ws.Range("H94").Value = 1090.1111
$ws.Range("I94").Value = 757.95654
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 757.95654
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -306.95654
$ws.Range("N94").Value = -3902

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 70000
$ws.Range("J14").Value = 70000
$ws.Range("L14").Value = 70000
$ws.Range("N14").Value = -70340
$ws.Range("H16").Value = 5918581.5
$ws.Range("I16").Value = 8548340
$ws.Range("J16").Value = 1625
$ws.Range("K16").Value = 8548340
$ws.Range("L16").Value = 1625
$ws.Range("M16").Value = -8548053
$ws.Range("N16").Value = -2199
$ws.Range("H31").Value = 7912.963
$ws.Range("I31").Value = 2459.3076
$ws.Range("J31").Value = 12977.071
$ws.Range("K31").Value = 2459.3076
$ws.Range("L31").Value = 12977.071
$ws.Range("M31").Value = -2164.3076
$ws.Range("N31").Value = -13567.071
$ws.Range("H34").Value = 7912.963
$ws.Range("I34").Value = 2459.3076
$ws.Range("J34").Value = 12977.071
$ws.Range("K34").Value = 2459.3076
$ws.Range("L34").Value = 12977.071
$ws.Range("M34").Value = -2257.3076
$ws.Range("N34").Value = -13381.071
$ws.Range("H42").Value = 6056
$ws.Range("I42").Value = 6056
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 6056
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -5463
$ws.Range("N42").Value = $null
$ws.Range("H94").Value = 3343.6667
$ws.Range("I94").Value = 3536.5557
$ws.Range("J94").Value = 3227.9333
$ws.Range("K94").Value = 3536.5557
$ws.Range("L94").Value = 3227.9333
$ws.Range("M94").Value = -3085.5557
$ws.Range("N94").Value = -4129.933300000001
$ws.Range("H113").Value = 5918581.5
$ws.Range("I113").Value = 8548340
$ws.Range("J113").Value = 1625
$ws.Range("K113").Value = 8548340
$ws.Range("L113").Value = 1625
$ws.Range("M113").Value = -8546170
$ws.Range("N113").Value = -5965
$ws.Range("H132").Value = 2689.1428
$ws.Range("I132").Value = 2479.818
$ws.Range("K132").Value = 7439.454000000001
$ws.Range("M132").Value = -4909.454000000001
$ws.Range("H134").Value = 2911.3
$ws.Range("I134").Value = 3390.5715
$ws.Range("J134").Value = 1793
$ws.Range("K134").Value = 10171.7145
$ws.Range("L134").Value = 5379
$ws.Range("M134").Value = -7636.7145
$ws.Range("N134").Value = -10449

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 300844.7
$ws.Range("I5").Value = 554.8889
$ws.Range("K5").Value = 1664.6667
$ws.Range("M5").Value = -1552.6667
$ws.Range("H107").Value = 476.84616
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 587.375
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1762.125
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -5602.125
$ws.Range("H113").Value = 2941817
$ws.Range("I113").Value = 6250481
$ws.Range("J113").Value = 782.2222
$ws.Range("K113").Value = 18751443
$ws.Range("L113").Value = 2346.6666
$ws.Range("M113").Value = -18749273
$ws.Range("N113").Value = -6686.6666
$ws.Range("H122").Value = 658.3333
$ws.Range("J122").Value = 787.7778
$ws.Range("L122").Value = 7090.000199999999
$ws.Range("N122").Value = -11990.0002
$ws.Range("H135").Value = 300844.7
$ws.Range("I135").Value = 554.8889
$ws.Range("K135").Value = 4994.0001
$ws.Range("M135").Value = -2459.0001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 90911040
$ws.Range("I40").Value = 125001800
$ws.Range("K40").Value = 125001800
$ws.Range("M40").Value = -125001664

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1654.7273
$ws.Range("I132").Value = 1026.0625
$ws.Range("J132").Value = 3331.1667
$ws.Range("K132").Value = 3078.1875
$ws.Range("L132").Value = 9993.500100000001
$ws.Range("M132").Value = -548.1875
$ws.Range("N132").Value = -15053.5001
